$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# Row 2 (Angelo): WIN_GATO 0 -> 10, new ULTIMO_JUEGO_GATO timestamp
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = "03/12/2022 15:08"

# Row 3 (Joseth): WIN_GATO 0 -> 2, new ULTIMO_JUEGO_GATO timestamp
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "03/12/2022 15:17"

# Row 4 (Jimena Hernandez): WIN_GATO 0 -> 1, new ULTIMO_JUEGO_GATO timestamp
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "03/12/2022 14:54"

# Row 5 (Isaac Mejias): WIN_GATO 0 -> 1, new ULTIMO_JUEGO_GATO timestamp,
# WIN_AHORCADO 0 -> 1, new ULTIMO_JUEGO_AHORCADO timestamp
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "03/12/2022 15:16"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "03/12/2022 14:37"
